$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 10561.565
$ws.Range("J17").Value = 10561.565
$ws.Range("L17").Value = 31684.695
$ws.Range("N17").Value = -32020.695
# Row 33
$ws.Range("H33").Value = 393.6111
$ws.Range("I33").Value = 227.25
$ws.Range("J33").Value = 1724.5
$ws.Range("K33").Value = 227.25
$ws.Range("L33").Value = 1724.5
$ws.Range("M33").Value = 1.75
$ws.Range("N33").Value = -2182.5
# Row 86
$ws.Range("H86").Value = 147731740
$ws.Range("I86").Value = 105268230
$ws.Range("J86").Value = 416667330
$ws.Range("K86").Value = 105268230
$ws.Range("L86").Value = 416667330
$ws.Range("M86").Value = -105267107
$ws.Range("N86").Value = -416669576
# Row 89
$ws.Range("H89").Value = 147731740
$ws.Range("I89").Value = 105268230
$ws.Range("J89").Value = 416667330
$ws.Range("K89").Value = 526341150
$ws.Range("L89").Value = 2083336650
$ws.Range("M89").Value = -526335534
$ws.Range("N89").Value = -2083347882
# Row 106
$ws.Range("H106").Value = 22228880
$ws.Range("I106").Value = 22228880
$ws.Range("K106").Value = 22228880
$ws.Range("M106").Value = -22228249
# Row 116
$ws.Range("H116").Value = 26329188
$ws.Range("I116").Value = 33348740
$ws.Range("J116").Value = 5868.5
$ws.Range("K116").Value = 33348740
$ws.Range("L116").Value = 5868.5
$ws.Range("M116").Value = -33345298
$ws.Range("N116").Value = -12752.5
# Row 137
$ws.Range("H137").Value = 1635965.4
$ws.Range("J137").Value = 3705930.2
$ws.Range("L137").Value = 11117790.6
$ws.Range("N137").Value = -11122890.6
# Row 138
$ws.Range("H138").Value = 2625.92
$ws.Range("I138").Value = 1427.2222
$ws.Range("J138").Value = 2744.4724
$ws.Range("K138").Value = 4281.6666
$ws.Range("L138").Value = 8233.4172
$ws.Range("M138").Value = 858.3334000000004
$ws.Range("N138").Value = -18513.4172

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 3412.8
$ws.Range("I74").Value = 2815.2856
$ws.Range("J74").Value = 6549.75
$ws.Range("K74").Value = 2815.2856
$ws.Range("L74").Value = 6549.75
$ws.Range("M74").Value = -1941.2856
$ws.Range("N74").Value = -8297.75
# Row 77
$ws.Range("H77").Value = 3412.8
$ws.Range("I77").Value = 2815.2856
$ws.Range("J77").Value = 6549.75
$ws.Range("K77").Value = 14076.428
$ws.Range("L77").Value = 32748.75
$ws.Range("M77").Value = -9708.428
$ws.Range("N77").Value = -41484.75

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 128500
$ws.Range("I20").Value = 251000
$ws.Range("K20").Value = 251000
$ws.Range("M20").Value = -250753
# Row 107
$ws.Range("H107").Value = 1494.3334
$ws.Range("I107").Value = 1362.5
$ws.Range("J107").Value = 1599.8
$ws.Range("K107").Value = 1362.5
$ws.Range("L107").Value = 1599.8
$ws.Range("M107").Value = 557.5
$ws.Range("N107").Value = -5439.8
# Row 134
$ws.Range("H134").Value = 1788781.4
$ws.Range("I134").Value = 2234127
$ws.Range("J134").Value = 7399.25
$ws.Range("K134").Value = 6702381
$ws.Range("L134").Value = 22197.75
$ws.Range("M134").Value = -6699846
$ws.Range("N134").Value = -27267.75
# Row 135
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 37.714287
$ws.Range("I7").Value = 34.916668
$ws.Range("K7").Value = 34.916668
$ws.Range("M7").Value = 78.083332
# Row 100
$ws.Range("H100").Value = 60005.89
$ws.Range("J100").Value = 60005.89
$ws.Range("L100").Value = 60005.89
$ws.Range("N100").Value = -62169.89
# Row 105
$ws.Range("H105").Value = 2613.8
$ws.Range("I105").Value = 1409.5
$ws.Range("J105").Value = 3416.6667
$ws.Range("K105").Value = 1409.5
$ws.Range("L105").Value = 3416.6667
$ws.Range("M105").Value = 337.5
$ws.Range("N105").Value = -6910.6667
# Row 114
$ws.Range("H114").Value = 58559.668
$ws.Range("J114").Value = 58559.668
$ws.Range("L114").Value = 58559.668
$ws.Range("N114").Value = -67237.66800000001
# Row 116
$ws.Range("H116").Value = 108753.5
$ws.Range("J116").Value = 108753.5
$ws.Range("L116").Value = 108753.5
$ws.Range("N116").Value = -117931.5
# Row 117
$ws.Range("H117").Value = 60666
$ws.Range("J117").Value = 60666
$ws.Range("L117").Value = 60666
$ws.Range("N117").Value = -69844
# Row 129
$ws.Range("H129").Value = 72568.5
$ws.Range("J129").Value = 86758
$ws.Range("L129").Value = 86758
$ws.Range("N129").Value = -96758
# Row 131
$ws.Range("H131").Value = 43886
$ws.Range("J131").Value = 43886
$ws.Range("L131").Value = 43886
$ws.Range("N131").Value = -53966
# Row 132
$ws.Range("H132").Value = 3923.6667
$ws.Range("I132").Value = 3914.8948
$ws.Range("K132").Value = 11744.6844
$ws.Range("M132").Value = -9214.6844

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 500
$ws.Range("J2").Value = 500
$ws.Range("L2").Value = 3000
$ws.Range("N2").Value = -3226
# Row 35
$ws.Range("H35").Value = 802
$ws.Range("I35").Value = 802
$ws.Range("K35").Value = 2406
$ws.Range("M35").Value = -2118
# Row 47
$ws.Range("H47").Value = 1352
$ws.Range("I47").Value = 700
$ws.Range("J47").Value = 2004
$ws.Range("K47").Value = 2100
$ws.Range("L47").Value = 6012
$ws.Range("M47").Value = -1669
$ws.Range("N47").Value = -6874
# Row 129
$ws.Range("H129").Value = 2253.375
$ws.Range("J129").Value = 2432.4285
$ws.Range("L129").Value = 7297.2855
$ws.Range("N129").Value = -17297.2855

$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 46026.57
$ws.Range("I52").Value = 46764.668
$ws.Range("K52").Value = 46764.668
$ws.Range("M52").Value = -46505.668

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2135.625
$ws.Range("I22").Value = 1748
$ws.Range("K22").Value = 1748
$ws.Range("M22").Value = -1453
# Row 27
$ws.Range("H27").Value = 2135.625
$ws.Range("I27").Value = 1748
$ws.Range("K27").Value = 1748
$ws.Range("M27").Value = -1641
# Row 58
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 100
$ws.Range("H100").Value = 3102
$ws.Range("I100").Value = 2630.25
$ws.Range("J100").Value = 4989
$ws.Range("K100").Value = 2630.25
$ws.Range("L100").Value = 4989
$ws.Range("M100").Value = -2089.25
$ws.Range("N100").Value = -6071
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# Row 132
$ws.Range("H132").Value = 3365.6667
$ws.Range("I132").Value = 3268.1538
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 9804.4614
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -7274.4614
$ws.Range("N132").Value = -17058.5

